$d = $word.ActiveDocument

# The title block originally had 4 short paragraphs:
#   1. "Isabel Finkbeiner, Vedant Nilabh, Ahmad Saeed"
#   2. "Professor Huskinson"
#   3. "CS 1340"
#   4. "Final Project: Report"
# The edit trims the author line down to a single name and removes the
# course/professor/title paragraphs entirely.

# Step 1: shorten the author line (paragraph 1) in place, preserving its
# run formatting (Times New Roman, 24 half-points, double spacing, etc.)
$found = $d.Paragraphs(1).Range.Find.Execute(
    "Isabel Finkbeiner, Vedant Nilabh, Ahmad Saeed",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Vedant Nilabh", 2)

# Step 2: delete paragraphs 2-4 ("Professor Huskinson", "CS 1340",
# "Final Project: Report") completely, paragraph marks included, so the
# document flows straight from the (now shortened) author line into the
# "Equity and Profitability..." title paragraph that used to be 5th.
$start = $d.Paragraphs(2).Range.Start
$end = $d.Paragraphs(4).Range.End
$d.Range($start, $end).Delete()
